$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmark Results")

$ws.Range("A2").Value = "4x4090"
$ws.Range("B2").Value = 1731.44
$ws.Range("C2").Value = 1.56
$ws.Range("D2").Value = 0.2502733755332748

$wb.Save()
